$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Fitness column (C) values per the diff
$ws.Range("C2").Value = 3833
$ws.Range("C3").Value = 3833
$ws.Range("C4").Value = 3833
$ws.Range("C5").Value = 4263
$ws.Range("C6").Value = 4263
$ws.Range("C7").Value = 4263
$ws.Range("C8").Value = 4263
$ws.Range("C9").Value = 4263
$ws.Range("C10").Value = 4518
$ws.Range("C11").Value = 4518
$ws.Range("C12").Value = 4531
